$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B width update ---
$ws.Columns("B").ColumnWidth = 14.6

# --- Update existing formula / data cells (rows 10-14) ---
# C10: new instance expression now builds a Y[] using the Y(...) constructor
$ws.Range("C10").Value = "'= new Y[] { new Y(`"a1`") } "

# C11: transform to x.name (was x.doubleValue())
$ws.Range("C11").Value = "'= `$S1[(X x) transform to x.name]"

# C12: transform unique to x.name
$ws.Range("C12").Value = "'= `$S1[(X x) transform unique to x.name]"

# C13: select all having x.name.length > 0
$ws.Range("C13").Value = "'= `$S1[(X x) select all having x.name.length > 0]"

# C14: split by x.name.length > 0
$ws.Range("C14").Value = "'= `$S1[(X x) split by x.name.length > 0]"

# --- New Datatype X block (rows 21-22) ---
$ws.Range("B21").Value = "Datatype X"
$ws.Range("B21:C21").Merge()
$ws.Range("B21:C21").HorizontalAlignment = -4108
$ws.Range("B21:C21").VerticalAlignment = -4108
$ws.Range("B21:C21").WrapText = $true

$ws.Range("B22").Value = "String"
$ws.Range("C22").Value = "name"
$ws.Range("B22:C22").HorizontalAlignment = -4108
$ws.Range("B22:C22").VerticalAlignment = -4108
$ws.Range("B22:C22").WrapText = $true

# --- New Datatype Y block (rows 25-26) ---
$ws.Range("B25").Value = "Datatype Y"
$ws.Range("B25:C25").Merge()
$ws.Range("B25:C25").HorizontalAlignment = -4108
$ws.Range("B25:C25").VerticalAlignment = -4108
$ws.Range("B25:C25").WrapText = $true

$ws.Range("B26").Value = "String"
$ws.Range("C26").Value = "name"
$ws.Range("B26:C26").HorizontalAlignment = -4108
$ws.Range("B26:C26").VerticalAlignment = -4108
$ws.Range("B26:C26").WrapText = $true

# --- Selection / view update ---
$ws.Range("F7").Select()
